# Add a new "2021" column (P) to the 4.2.2 indicator table, mirroring the
# existing 2020 column (O) — same formatting, new year label + value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column O into a freshly inserted column P so the new column
# inherits O's exact cell styles (header style on row 4, value style on
# row 5) instead of Excel computing a brand-new style from scratch.
$ws.Columns("O").Copy()
$ws.Columns("P").Insert(-4161, 0)

# Overwrite the copied (2020) values with the new 2021 figures.
$ws.Range("P4").Value = 2021
$ws.Range("P5").Value = 80.9

# Match the author's final selection in the saved file.
$ws.Range("N10").Select()
